$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; B = 4.633760690483867; C = 0.9157948689166304 },
    @{ Row = 3; B = 5.099333848387586; C = 1.810025689074652 },
    @{ Row = 4; B = 5.194255551181947; C = 2.841746928840663 },
    @{ Row = 5; B = 9.760021383178435; C = 3.833084125498456 },
    @{ Row = 6; B = 14.12308159874659; C = 4.754056727903757 },
    @{ Row = 7; B = 16.52380481384802; C = 5.821134536603146 },
    @{ Row = 8; B = 16.82933328193898; C = 6.823351705313178 },
    @{ Row = 9; B = 20.03966007549048; C = 7.928030773264516 },
    @{ Row = 10; B = 23.72023758006711; C = 9.536905490795473 },
    @{ Row = 11; B = 33.61922367188829; C = 10.43877629938486 },
    @{ Row = 12; B = 37.47745075236205; C = 11.35100749197012 },
    @{ Row = 13; B = 37.55693407908596; C = 12.3259755717076 },
    @{ Row = 14; B = 37.83548351054211; C = 13.64655386399941 },
    @{ Row = 15; B = 40.89463314041754; C = 14.51753660139681 },
    @{ Row = 16; B = 42.26911647374695; C = 15.68379542862329 },
    @{ Row = 17; B = 43.51841339956567; C = 16.54065194121847 },
    @{ Row = 18; B = 45.74305845682478; C = 17.75134382239723 },
    @{ Row = 19; B = 48.48031580598457; C = 18.62451020148014 },
    @{ Row = 20; B = 48.70006017335604; C = 19.63965975641869 },
    @{ Row = 21; B = 50.72028802835469; C = 20.68283995594952 },
    @{ Row = 22; B = 53.15185544257069; C = 21.57835024785508 },
    @{ Row = 23; B = 53.43996451828138; C = 22.53708503187125 },
    @{ Row = 24; B = 53.51230141012838; C = 23.50006041030797 },
    @{ Row = 25; B = 63.01920475955566; C = 24.47509269906413 },
    @{ Row = 26; B = 63.08717969310466; C = 25.47784677251192 },
    @{ Row = 27; B = 63.64371229842498; C = 26.37499013797765 },
    @{ Row = 28; B = 64.47967849202442; C = 27.46880624219978 },
    @{ Row = 29; B = 68.12542564799607; C = 28.39416604254507 },
    @{ Row = 30; B = 68.2153168382511; C = 29.4573819646596 },
    @{ Row = 31; B = 69.16480061400163; C = 30.74857370857071 },
    @{ Row = 32; B = 69.72143464578903; C = 31.59851834108667 },
    @{ Row = 33; B = 74.07273048436051; C = 32.47609756478946 },
    @{ Row = 34; B = 76.76810005691152; C = 33.42884756270491 },
    @{ Row = 35; B = 78.00203709683637; C = 34.33068062681019 },
    @{ Row = 36; B = 78.3469708474522; C = 35.31542525016336 },
    @{ Row = 37; B = 81.61692436282384; C = 36.34145238766495 },
    @{ Row = 38; B = 85.87418795504234; C = 37.40668144375827 },
    @{ Row = 39; B = 87.46854372792866; C = 38.35944274176383 },
    @{ Row = 40; B = 87.54137630355916; C = 39.45126764029932 },
    @{ Row = 41; B = 89.45817422194929; C = 40.53636155876845 },
    @{ Row = 42; B = 90.07811632152834; C = 41.49976361373444 },
    @{ Row = 43; B = 93.57312516855718; C = 43.05420825428281 },
    @{ Row = 44; B = 95.45341764459619; C = 44.00258206439221 },
    @{ Row = 45; B = 97.01825787982882; C = 45.05893695464815 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}

$ws.Rows("46:49").Delete()
